$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 13:01"

# Update country stats rows that changed (re-sorted by Casos totales desc, with refreshed data)
$ws.Cells.Item(13, 1).Value = "Iran"
$ws.Cells.Item(13, 2).Value = 232863
$ws.Cells.Item(13, 3).Value = 2652
$ws.Cells.Item(13, 4).Value = 194098
$ws.Cells.Item(13, 5).Value = 27659
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 148
$ws.Cells.Item(13, 8).Value = 11106

$ws.Cells.Item(14, 1).Value = "Mexico"
$ws.Cells.Item(14, 2).Value = 231770
$ws.Cells.Item(14, 3).Value = 5681
$ws.Cells.Item(14, 4).Value = 138319
$ws.Cells.Item(14, 5).Value = 64941
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 741
$ws.Cells.Item(14, 8).Value = 28510

$ws.Cells.Item(17, 1).Value = "Alemania"
$ws.Cells.Item(17, 2).Value = 196361
$ws.Cells.Item(17, 3).Value = 37
$ws.Cells.Item(17, 4).Value = 180300
$ws.Cells.Item(17, 5).Value = 7000
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 9061

$ws.Cells.Item(29, 1).Value = "Bielorrusia"
$ws.Cells.Item(29, 2).Value = 62698
$ws.Cells.Item(29, 3).Value = 274
$ws.Cells.Item(29, 4).Value = 48738
$ws.Cells.Item(29, 5).Value = 13555
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 7
$ws.Cells.Item(29, 8).Value = 405

$ws.Cells.Item(48, 1).Value = "Suiza"
$ws.Cells.Item(48, 2).Value = 31967
$ws.Cells.Item(48, 3).Value = 116
$ws.Cells.Item(48, 4).Value = 29200
$ws.Cells.Item(48, 5).Value = 802
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 1965

$ws.Cells.Item(49, 1).Value = "Rumania"
$ws.Cells.Item(49, 2).Value = 27746
$ws.Cells.Item(49, 3).Value = 450
$ws.Cells.Item(49, 4).Value = 19363
$ws.Cells.Item(49, 5).Value = 6696
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 20
$ws.Cells.Item(49, 8).Value = 1687

$ws.Cells.Item(50, 1).Value = "Barein"
$ws.Cells.Item(50, 2).Value = 27414
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 21948
$ws.Cells.Item(50, 5).Value = 5373
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 93

$ws.Cells.Item(63, 1).Value = "Nepal"
$ws.Cells.Item(63, 2).Value = 14519
$ws.Cells.Item(63, 3).Value = 473
$ws.Cells.Item(63, 4).Value = 5320
$ws.Cells.Item(63, 5).Value = 9168
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 31

$ws.Cells.Item(64, 1).Value = "Argelia"
$ws.Cells.Item(64, 2).Value = 14272
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 4).Value = 10040
$ws.Cells.Item(64, 5).Value = 3312
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 920

$ws.Cells.Item(77, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(77, 2).Value = 7189
$ws.Cells.Item(77, 3).Value = 67
$ws.Cells.Item(77, 4).Value = 2317
$ws.Cells.Item(77, 5).Value = 4696
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 176

$ws.Cells.Item(78, 1).Value = "Senegal"
$ws.Cells.Item(78, 2).Value = 7054
$ws.Cells.Item(78, 3).Value = 129
$ws.Cells.Item(78, 4).Value = 4599
$ws.Cells.Item(78, 5).Value = 2334
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 5
$ws.Cells.Item(78, 8).Value = 121

$ws.Cells.Item(79, 1).Value = "El Salvador"
$ws.Cells.Item(79, 2).Value = 7000
$ws.Cells.Item(79, 3).Value = 264
$ws.Cells.Item(79, 4).Value = 4115
$ws.Cells.Item(79, 5).Value = 2694
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 9
$ws.Cells.Item(79, 8).Value = 191

$ws.Cells.Item(90, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(90, 2).Value = 4788
$ws.Cells.Item(90, 3).Value = 182
$ws.Cells.Item(90, 4).Value = 2515
$ws.Cells.Item(90, 5).Value = 2084
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 189

$ws.Cells.Item(91, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(91, 2).Value = 4704
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 4550
$ws.Cells.Item(91, 5).Value = 99
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 55

$ws.Cells.Item(100, 1).Value = "Estado de Palestina"
$ws.Cells.Item(100, 2).Value = 2978
$ws.Cells.Item(100, 3).Value = 220
$ws.Cells.Item(100, 4).Value = 460
$ws.Cells.Item(100, 5).Value = 2510
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 8

$ws.Cells.Item(101, 1).Value = "Somalia"
$ws.Cells.Item(101, 2).Value = 2924
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 932
$ws.Cells.Item(101, 5).Value = 1902
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 90

$ws.Cells.Item(102, 1).Value = "Croacia"
$ws.Cells.Item(102, 2).Value = 2831
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 2155
$ws.Cells.Item(102, 5).Value = 568
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 108

$ws.Cells.Item(106, 1).Value = "Madagascar"
$ws.Cells.Item(106, 2).Value = 2403
$ws.Cells.Item(106, 3).Value = 100
$ws.Cells.Item(106, 4).Value = 1040
$ws.Cells.Item(106, 5).Value = 1339
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 2
$ws.Cells.Item(106, 8).Value = 24

$ws.Cells.Item(107, 1).Value = "Maldivas"
$ws.Cells.Item(107, 2).Value = 2382
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 1954
$ws.Cells.Item(107, 5).Value = 419
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 9

$ws.Cells.Item(108, 1).Value = "Cuba"
$ws.Cells.Item(108, 2).Value = 2348
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 2218
$ws.Cells.Item(108, 5).Value = 44
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 86

$ws.Cells.Item(114, 1).Value = "Estonia"
$ws.Cells.Item(114, 2).Value = 1990
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = 1842
$ws.Cells.Item(114, 5).Value = 79
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 69

$ws.Cells.Item(140, 1).Value = "Uganda"
$ws.Cells.Item(140, 2).Value = 902
$ws.Cells.Item(140, 3).Value = 9
$ws.Cells.Item(140, 4).Value = 847
$ws.Cells.Item(140, 5).Value = 55
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 0

$ws.Cells.Item(150, 1).Value = "Malta"
$ws.Cells.Item(150, 2).Value = 671
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 649
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 9

$ws.Cells.Item(161, 1).Value = "Siria"
$ws.Cells.Item(161, 2).Value = 312
$ws.Cells.Item(161, 3).Value = 19
$ws.Cells.Item(161, 4).Value = 113
$ws.Cells.Item(161, 5).Value = 190
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 9

$ws.Cells.Item(162, 1).Value = "Birmania"
$ws.Cells.Item(162, 2).Value = 304
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(162, 4).Value = 222
$ws.Cells.Item(162, 5).Value = 76
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 6

$ws.Cells.Item(163, 1).Value = "Comoras"
$ws.Cells.Item(163, 2).Value = 303
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 200
$ws.Cells.Item(163, 5).Value = 96
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 7

$ws.Cells.Item(164, 1).Value = "Namibia"
$ws.Cells.Item(164, 2).Value = 293
$ws.Cells.Item(164, 3).Value = 8
$ws.Cells.Item(164, 4).Value = 24
$ws.Cells.Item(164, 5).Value = 269
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 0

$ws.Cells.Item(205, 1).Value = "Dominica"
$ws.Cells.Item(205, 2).Value = 18
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 18
$ws.Cells.Item(205, 5).Value = 0
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

$ws.Cells.Item(206, 1).Value = "Fiyi"
$ws.Cells.Item(206, 2).Value = 18
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 18
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0
